# Scheduled price-data refresh: update computed profit columns (H:N)
# on each affected Leve row across sheets, per the upstream market-board sync.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18: You Grow, Girl | Growth Formula Beta
$ws.Range("H18").Value = 1096.1111
$ws.Range("I18").Value = 782
$ws.Range("J18").Value = 2666.6667
$ws.Range("K18").Value = 782
$ws.Range("L18").Value = 2666.6667
$ws.Range("M18").Value = -498
$ws.Range("N18").Value = -3234.6667

# Row 64: Forged from the Void | Void Glue
$ws.Range("H64").Value = 86232.414
$ws.Range("J64").Value = 3298.1667
$ws.Range("L64").Value = 3298.1667
$ws.Range("N64").Value = -3794.1667

# Row 67: Dodging the Draft (L) | Void Glue
$ws.Range("H67").Value = 86232.414
$ws.Range("J67").Value = 3298.1667
$ws.Range("L67").Value = 3298.1667
$ws.Range("N67").Value = -5014.1667

# Row 129: Practical Command | Commanding Craftsman's Draught
$ws.Range("H129").Value = 3023.8044
$ws.Range("I129").Value = 14844.571
$ws.Range("J129").Value = 902.12823
$ws.Range("K129").Value = 44533.713
$ws.Range("L129").Value = 2706.38469
$ws.Range("M129").Value = -39533.713
$ws.Range("N129").Value = -12706.38469

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 3681276
$ws.Range("I132").Value = 4036980
$ws.Range("K132").Value = 12110940
$ws.Range("M132").Value = -12108410

# Row 135: For Tired Minds | Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 4739.6665
$ws.Range("I135").Value = 1770.75
$ws.Range("J135").Value = 5989.737
$ws.Range("K135").Value = 15936.75
$ws.Range("L135").Value = 53907.633
$ws.Range("M135").Value = -13401.75
$ws.Range("N135").Value = -58977.633

# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 2046.25
$ws.Range("I137").Value = 1738.6364
$ws.Range("K137").Value = 5215.9092
$ws.Range("M137").Value = -2665.9092

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 23851.217
$ws.Range("I32").Value = 7366.643
$ws.Range("J32").Value = 196939.25
$ws.Range("K32").Value = 7366.643
$ws.Range("L32").Value = 196939.25
$ws.Range("M32").Value = -7079.643
$ws.Range("N32").Value = -197513.25

# Row 46: Get Me the Usual | Heavy Steel Flanchard
$ws.Range("H46").Value = 2874.5715
$ws.Range("J46").Value = 2520.3333
$ws.Range("L46").Value = 2520.3333
$ws.Range("N46").Value = -3158.3333

# Row 86: Sir, Dost Thou Even Heft | Adamantite Chain Hose of Fending
$ws.Range("H86").Value = 31854.285
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 31854.285
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 31854.285
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -34226.285

# Row 89: Men in Adamantite (L) | Adamantite Chain Hose of Fending
$ws.Range("H89").Value = 31854.285
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 31854.285
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 95562.855
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -107418.855

# Row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 1992
$ws.Range("I122").Value = 1681.6086
$ws.Range("J122").Value = 2641
$ws.Range("K122").Value = 5044.825800000001
$ws.Range("L122").Value = 7923
$ws.Range("M122").Value = -2594.825800000001
$ws.Range("N122").Value = -12823

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 13017.48
$ws.Range("I132").Value = 16099.8
$ws.Range("J132").Value = 2743.0833
$ws.Range("K132").Value = 48299.39999999999
$ws.Range("L132").Value = 8229.249899999999
$ws.Range("M132").Value = -45769.39999999999
$ws.Range("N132").Value = -13289.2499

$ws = $wb.Worksheets.Item("BSM")
# Row 36: I Saw What You Did There | Iron Chocobotail Saw
$ws.Range("H36").Value = 12260.167
$ws.Range("I36").Value = 12260.167
$ws.Range("K36").Value = 12260.167
$ws.Range("M36").Value = -11726.167

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 35918.617
$ws.Range("I31").Value = 905.8333
$ws.Range("J31").Value = 82602.336
$ws.Range("K31").Value = 905.8333
$ws.Range("L31").Value = 82602.336
$ws.Range("M31").Value = -610.8333
$ws.Range("N31").Value = -83192.336

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 35918.617
$ws.Range("I34").Value = 905.8333
$ws.Range("J34").Value = 82602.336
$ws.Range("K34").Value = 905.8333
$ws.Range("L34").Value = 82602.336
$ws.Range("M34").Value = -703.8333
$ws.Range("N34").Value = -83006.336

# Row 125: A Wristy Experiment | Palm Bracelets of Aiming
$ws.Range("H125").Value = 30000
$ws.Range("J125").Value = 30000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -34920

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 2635.2068
$ws.Range("I132").Value = 2629.0417
$ws.Range("K132").Value = 7887.125100000001
$ws.Range("M132").Value = -5357.125100000001

$ws = $wb.Worksheets.Item("CUL")
# Row 41: Gegeruju Gets Down | Cornbread
$ws.Range("H41").Value = 2450
$ws.Range("J41").Value = 4400
$ws.Range("L41").Value = 13200
$ws.Range("N41").Value = -13876

# Row 82: Persuasion of a Higher Power | Baked Pipira Pira
$ws.Range("H82").Value = 2419.9092
$ws.Range("I82").Value = 2088.4285
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 6265.2855
$ws.Range("L82").Value = 9000
$ws.Range("M82").Value = -5859.2855
$ws.Range("N82").Value = -9812

# Row 85: Loaves and Fishes (L) | Baked Pipira Pira
$ws.Range("H85").Value = 2419.9092
$ws.Range("I85").Value = 2088.4285
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 6265.2855
$ws.Range("L85").Value = 9000
$ws.Range("M85").Value = -4861.2855
$ws.Range("N85").Value = -11808

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit | Mythrite Ingot
$ws.Range("H70").Value = 80552.7
$ws.Range("I70").Value = 116827.61
$ws.Range("K70").Value = 116827.61
$ws.Range("M70").Value = -116557.61

# Row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
$ws.Range("H73").Value = 80552.7
$ws.Range("I73").Value = 116827.61
$ws.Range("K73").Value = 116827.61
$ws.Range("M73").Value = -115891.61

# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 1806.4
$ws.Range("I122").Value = 1111.6666
$ws.Range("J122").Value = 2848.5
$ws.Range("K122").Value = 3334.9998
$ws.Range("L122").Value = 8545.5
$ws.Range("M122").Value = -884.9998000000001
$ws.Range("N122").Value = -13445.5

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 2571.1707
$ws.Range("I132").Value = 1791.7941
$ws.Range("K132").Value = 5375.3823
$ws.Range("M132").Value = -2845.3823

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad | Toad Leather
$ws.Range("H40").Value = 55640.95
$ws.Range("I40").Value = 128091
$ws.Range("K40").Value = 128091
$ws.Range("M40").Value = -127955

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 3121.5151
$ws.Range("I132").Value = 3152.68
$ws.Range("K132").Value = 9458.039999999999
$ws.Range("M132").Value = -6928.039999999999

# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("H136").Value = 1807.5349
$ws.Range("I136").Value = 1241.6666
$ws.Range("J136").Value = 3674.9
$ws.Range("K136").Value = 3724.9998
$ws.Range("L136").Value = 11024.7
$ws.Range("M136").Value = -1174.9998
$ws.Range("N136").Value = -16124.7

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches | Crawler Silk
$ws.Range("H81").Value = 288482.56
$ws.Range("I81").Value = 334600
$ws.Range("J81").Value = 253894.5
$ws.Range("K81").Value = 669200
$ws.Range("L81").Value = 507789
$ws.Range("M81").Value = -668139
$ws.Range("N81").Value = -509911

# Row 84: To Kill a Dragon on Nameday (L) | Crawler Silk
$ws.Range("H84").Value = 288482.56
$ws.Range("I84").Value = 334600
$ws.Range("J84").Value = 253894.5
$ws.Range("K84").Value = 3346000
$ws.Range("L84").Value = 2538945
$ws.Range("M84").Value = -3340696
$ws.Range("N84").Value = -2549553

# Row 115: Gloves Come in Handy | Pixie Cotton Sleeves of Crafting
$ws.Range("H115").Value = 38443.445
$ws.Range("J115").Value = 38443.445
$ws.Range("L115").Value = 38443.445
$ws.Range("N115").Value = -41577.445

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 4893.976
$ws.Range("I132").Value = 2573.9656
$ws.Range("J132").Value = 10069.385
$ws.Range("K132").Value = 7721.8968
$ws.Range("L132").Value = 30208.155
$ws.Range("M132").Value = -5191.8968
$ws.Range("N132").Value = -35268.155
